$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 6 values (new unique shared strings)
$ws.Range("A6").Value = "Source4"
$ws.Range("B6").Value = "Drop-down"
$ws.Range("C6").Value = "CapDentalBaseClaimData.source4"
$ws.Range("D6").Value = "Rule1: `nRule description: Field is mandatory `nRuleName:Source4Mandatory `nError Message: Source4is required."

# Match formatting used by the other "Drop-down" rows (row 3 pattern)
$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Size = 9
$ws.Range("B6").Font.Color = 4668466

$ws.Range("C6").Font.Name = "Times New Roman"
$ws.Range("C6").Font.Size = 10
$ws.Range("C6").Font.Bold = $true

$ws.Range("D6").Font.Name = "Calibri"
$ws.Range("D6").Font.Size = 11
$ws.Range("D6").WrapText = $true

$ws.Rows.Item(6).RowHeight = 159

$ws.Range("D6").Select()
